# Update "想去人数" (number of people who want to go) counts on the
# 展览 (Exhibitions) sheet and the 全部类型 (All Types) sheet.
# Each row is matched by event name in column C so the correct row is
# updated regardless of its position on each sheet.

$wb = $excel.ActiveWorkbook

$updates = @{
    "信州·上漫·ACG动漫游戏嘉年华"                 = 60
    "南昌·鸢歌弦 代号鸢同人only"                   = 159
    "上饶·星河城市动漫文化节"                      = 355
    "南昌·CM04动漫游戏博览会"                      = 5256
    "南昌·云芽动漫音乐嘉年华"                      = 5340
    "南昌·云芽动漫音乐嘉年华·封茗囧菌内场票"        = 618
    "上饶·ETI03动漫节"                             = 5
    "南昌·萌卡动漫展"                              = 1364
    "九江·第二届异次元动漫嘉年华"                  = 104
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $name = $ws.Cells.Item($r, 3).Value2
        if ($null -ne $name -and $updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value2 = $updates[$name]
        }
    }
}
